$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Delete row 5 (event removed)
$ws.Rows.Item(5).Delete()

# Row 2
$ws.Range("A2").Value = 1369
$ws.Range("B2").Value = '2025-11-24T19:00:00'
$ws.Range("C2").Value = 'Локомотив'
$ws.Range("D2").Value = 'Сибирь'
$ws.Range("E2").Value = 897811
$ws.Range("F2").Value = 'https://text.khl.ru/text/897811.html'
$ws.Range("G2").Value = 3.816985
$ws.Range("H2").Value = 1.117647
$ws.Range("I2").Value = 2.877774
$ws.Range("J2").Value = 5.647059
$ws.Range("K2").Value = 4.732022
$ws.Range("L2").Value = 1.997711
$ws.Range("M2").Value = 4.934632
$ws.Range("N2").Value = 35.165054
$ws.Range("O2").Value = 20.258009
$ws.Range("P2").Value = 55.423063
$ws.Range("Q2").Value = 0.067386
$ws.Range("R2").Value = -0.2
$ws.Range("S2").Value = 0.79889
$ws.Range("T2").Value = 0.092067
$ws.Range("U2").Value = 0.099583
$ws.Range("V2").Value = 0.096988
$ws.Range("W2").Value = 0.893551
$ws.Range("X2").Value = 0.199104
$ws.Range("Y2").Value = 0.791435
$ws.Range("Z2").Value = 0.336547
$ws.Range("AA2").Value = 0.653992
$ws.Range("AB2").Value = 0.490706
$ws.Range("AC2").Value = 0.499833
$ws.Range("AD2").Value = 0.638913
$ws.Range("AE2").Value = 0.351627
$ws.Range("AF2").Value = 0.949509
$ws.Range("AG2").Value = 0.050491
$ws.Range("AH2").Value = 0.850887
$ws.Range("AI2").Value = 0.149113
$ws.Range("AJ2").Value = 0.593374
$ws.Range("AK2").Value = 0.406626
$ws.Range("AL2").Value = 0.322704
$ws.Range("AM2").Value = 0.677296
$ws.Range("AN2").Value = 0.945669
$ws.Range("AO2").Value = 0.321248

# Row 3
$ws.Range("A3").Value = 1369
$ws.Range("B3").Value = '2025-11-24T19:30:00'
$ws.Range("C3").Value = 'Динамо М'
$ws.Range("D3").Value = 'Амур'
$ws.Range("E3").Value = 897809
$ws.Range("F3").Value = 'https://text.khl.ru/text/897809.html'
$ws.Range("G3").Value = 1.983158
$ws.Range("H3").Value = 3.411544
$ws.Range("I3").Value = 3.07998
$ws.Range("J3").Value = 4.220485
$ws.Range("K3").Value = 3.101821
$ws.Range("L3").Value = 3.245762
$ws.Range("M3").Value = 5.394701
$ws.Range("N3").Value = 29.463809
$ws.Range("O3").Value = 33.087107
$ws.Range("P3").Value = 62.550916
$ws.Range("Q3").Value = -0.105359
$ws.Range("R3").Value = 0.178929
$ws.Range("S3").Value = 0.396516
$ws.Range("T3").Value = 0.161561
$ws.Range("U3").Value = 0.44098
$ws.Range("V3").Value = 0.122777
$ws.Range("W3").Value = 0.87628
$ws.Range("X3").Value = 0.241218
$ws.Range("Y3").Value = 0.757839
$ws.Range("Z3").Value = 0.391581
$ws.Range("AA3").Value = 0.607476
$ws.Range("AB3").Value = 0.550654
$ws.Range("AC3").Value = 0.448403
$ws.Range("AD3").Value = 0.694902
$ws.Range("AE3").Value = 0.304155
$ws.Range("AF3").Value = 0.815552
$ws.Range("AG3").Value = 0.184448
$ws.Range("AH3").Value = 0.599231
$ws.Range("AI3").Value = 0.400769
$ws.Range("AJ3").Value = 0.834675
$ws.Range("AK3").Value = 0.165325
$ws.Range("AL3").Value = 0.629565
$ws.Range("AM3").Value = 0.370435
$ws.Range("AN3").Value = 0.709696
$ws.Range("AO3").Value = 0.747436

# Row 4
$ws.Range("A4").Value = 1369
$ws.Range("B4").Value = '2025-11-24T19:30:00'
$ws.Range("C4").Value = 'ЦСКА'
$ws.Range("D4").Value = 'СКА'
$ws.Range("E4").Value = 897810
$ws.Range("F4").Value = 'https://text.khl.ru/text/897810.html'
$ws.Range("G4").Value = 4.615385
$ws.Range("H4").Value = 2.109924
$ws.Range("I4").Value = 1.115385
$ws.Range("J4").Value = 1.63253
$ws.Range("K4").Value = 3.123958
$ws.Range("L4").Value = 1.612654
$ws.Range("M4").Value = 6.725309
$ws.Range("N4").Value = 34.357697
$ws.Range("O4").Value = 28.798769
$ws.Range("P4").Value = 63.156466
$ws.Range("Q4").Value = 0.2
$ws.Range("R4").Value = -0.013181
$ws.Range("S4").Value = 0.676504
$ws.Range("T4").Value = 0.151808
$ws.Range("U4").Value = 0.171278
$ws.Range("V4").Value = 0.30396
$ws.Range("W4").Value = 0.695631
$ws.Range("X4").Value = 0.487857
$ws.Range("Y4").Value = 0.511733
$ws.Range("Z4").Value = 0.662067
$ws.Range("AA4").Value = 0.337523
$ws.Range("AB4").Value = 0.799595
$ws.Range("AC4").Value = 0.199996
$ws.Range("AD4").Value = 0.892654
$ws.Range("AE4").Value = 0.106936
$ws.Range("AF4").Value = 0.818617
$ws.Range("AG4").Value = 0.181383
$ws.Range("AH4").Value = 0.604001
$ws.Range("AI4").Value = 0.395999
$ws.Range("AJ4").Value = 0.479147
$ws.Range("AK4").Value = 0.520853
$ws.Range("AL4").Value = 0.219917
$ws.Range("AM4").Value = 0.780083
$ws.Range("AN4").Value = 0.924297
$ws.Range("AO4").Value = 0.509023
